# Update countries & provincias Spain
# Applies the May 1 2020 20:52 -> 21:22 data refresh:
#  - Updated case counts for Estados Unidos (row 4), Alemania (row 9), Costa Rica (row 100)
#  - Updated case counts for Yemen, which overtakes Butan in "Casos totales" and swaps
#    places with it in the sorted country list (rows 210/211)
#  - Updated the "Datos actualizados" timestamp banner

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1115977
$ws.Range("C4").Value = 20954
$ws.Range("D4").Value = 158287
$ws.Range("E4").Value = 892661
$ws.Range("F4").Value = 15151
$ws.Range("G4").Value = 1173
$ws.Range("H4").Value = 65029

# --- Row 9: Alemania ---
$ws.Range("B9").Value = 163759
$ws.Range("C9").Value = 750
$ws.Range("D9").Value = 126900
$ws.Range("E9").Value = 30197
$ws.Range("F9").Value = 2189
$ws.Range("G9").Value = 39
$ws.Range("H9").Value = 6662

# --- Row 100: Costa Rica ---
$ws.Range("B100").Value = 725
$ws.Range("C100").Value = 6
$ws.Range("D100").Value = 355
$ws.Range("E100").Value = 364
$ws.Range("F100").Value = 7
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 6

# --- Rows 210/211: Yemen's numbers increase to 7 total cases, tying then
#     overtaking Butan in the descending sort, so they swap row positions ---
$ws.Range("A210").Value = "Yemen"
$ws.Range("B210").Value = 7
$ws.Range("C210").Value = 1
$ws.Range("D210").Value = 1
$ws.Range("E210").Value = 4
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 2

$ws.Range("A211").Value = "Butan"
$ws.Range("B211").Value = 7
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 5
$ws.Range("E211").Value = 2
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# --- Timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 21:22"
